$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add the new question column (J) -- question5 / dropdown question
$ws.Cells.Item(1, 10).Value = "question5"
$ws.Cells.Item(2, 10).Value = "dropdown|What is your favorite movie?|Anything Monty Python|Something other than Monty Python"

# Update the active selection to match the post-edit state
$ws.Range("H10").Select() | Out-Null
